$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Quarterly")

# Insert a new column before the existing "P/l before tax" column (old L), shifting
# old columns L:T one position to the right (M:U). The new column will hold the
# "Exceptional items" figures that sit between "P/l before exceptional items & tax"
# (col K) and "P/l before tax" (now col M).
$ws.Range("L1").EntireColumn.Insert()

# Header row (raw labels, row 1) and title-case header row (row 2) for the new column.
$ws.Range("L1").Value = "Exceptional items"
$ws.Range("L2").Value = "Exceptional Items"

# Populate the new "Exceptional items" column wherever P/l before tax (now col M)
# differs from P/l before exceptional items & tax (col K) - i.e. quarters that
# actually booked an exceptional item.
$rows = @(37, 39, 40, 41, 44, 46)
foreach ($r in $rows) {
    $k = $ws.Range("K$r").Value2
    $m = $ws.Range("M$r").Value2
    $ws.Range("L$r").Value = $m - $k
}
